# Add columns I (I0) and J (IF) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 1), styled like the existing headers in column H.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-21.
$dataI = @(6, 7, 6, 7, 7, 5, 6, 6, 3, 6, 9, 9, 9, 9, 9, 3, 7, 7, 8, 5)
$dataJ = @(8, 8, 8, 8, 8, 7, 7, 7, 5, 7, 9, 9, 9, 9, 9, 3, 7, 7, 8, 5)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
